$d = $word.ActiveDocument

# The "accuarcy results before removing outliers" table is the 2nd table
# in the document (index 2). Its 5th column (SVM) grows from 1806 dxa to
# 1866 dxa (widening the whole table from 9026 to 9086 dxa), and the
# previously-empty "Validation Accuracy" / "Test Accuracy" rows get filled
# in with the actual computed accuracy numbers now that the "remove
# outliers" step (and its separate before/after comparison) was dropped.
$t = $d.Tables.Item(2)

# Widen the SVM column (updates tblGrid + every row's tcW together).
$t.Columns.Item(5).Width = 93.3
# Keep the table's overall preferred width in sync (9086 dxa = 454.3 pt).
$t.PreferredWidth = 454.3

# Row 2 = "Validation Accuracy"
$t.Cell(2, 2).Range.Text = "0.686652391"
$t.Cell(2, 3).Range.Text = "0.7323340471"
$t.Cell(2, 4).Range.Text = "0.729478943611706"
$t.Cell(2, 5).Range.Text = "0.7323340471"

# Row 3 = "Test Accuracy"
$t.Cell(3, 2).Range.Text = "0.6842105263"
$t.Cell(3, 3).Range.Text = "0.7343117408"
$t.Cell(3, 4).Range.Text = "0.7211538461538461"
$t.Cell(3, 5).Range.Text = "0.7277327935"
